$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 18 (item id 5471)
$ws.Range("H18").Value = 2298.8
$ws.Range("I18").Value = 2311
$ws.Range("K18").Value = 2311
$ws.Range("M18").Value = -2027
# Row 19 (item id 7015)
$ws.Range("H19").Value = 882.9
$ws.Range("I19").Value = 328.8889
$ws.Range("J19").Value = 1336.1818
$ws.Range("K19").Value = 328.8889
$ws.Range("L19").Value = 1336.1818
$ws.Range("M19").Value = -153.8889
$ws.Range("N19").Value = -1686.1818
# Row 62 (item id 27781)
$ws.Range("H62").Value = 3187.4
$ws.Range("I62").Value = 2829.5
$ws.Range("J62").Value = 3724.25
$ws.Range("K62").Value = 2829.5
$ws.Range("L62").Value = 3724.25
$ws.Range("M62").Value = -2205.5
$ws.Range("N62").Value = -4972.25
# Row 65 (item id 27781)
$ws.Range("H65").Value = 3187.4
$ws.Range("I65").Value = 2829.5
$ws.Range("J65").Value = 3724.25
$ws.Range("K65").Value = 14147.5
$ws.Range("L65").Value = 18621.25
$ws.Range("M65").Value = -11027.5
$ws.Range("N65").Value = -24861.25
# Row 116 (item id 27778)
$ws.Range("H116").Value = 2612.1785
$ws.Range("I116").Value = 2065.6316
$ws.Range("K116").Value = 2065.6316
$ws.Range("M116").Value = 1376.3684
# Row 132 (item id 44049)
$ws.Range("H132").Value = 9262066
$ws.Range("I132").Value = 12823500
$ws.Range("J132").Value = 2337.6
$ws.Range("K132").Value = 38470500
$ws.Range("L132").Value = 7012.799999999999
$ws.Range("M132").Value = -38467970
$ws.Range("N132").Value = -12072.8
# Row 137 (item id 44013)
$ws.Range("H137").Value = 1901.625
$ws.Range("J137").Value = 2650
$ws.Range("L137").Value = 7950
$ws.Range("N137").Value = -13050

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61 (item id 43999)
$ws.Range("H61").Value = 901.4583
$ws.Range("I61").Value = 901.4583
$ws.Range("K61").Value = 901.4583
$ws.Range("M61").Value = -689.4583
# Row 122 (item id 36168)
$ws.Range("H122").Value = 1453.0667
$ws.Range("I122").Value = 1369
$ws.Range("K122").Value = 4107
$ws.Range("M122").Value = -1657
# Row 132 (item id 43997)
$ws.Range("H132").Value = 2925.6453
$ws.Range("I132").Value = 2769.72
$ws.Range("J132").Value = 3575.3333
$ws.Range("K132").Value = 8309.16
$ws.Range("L132").Value = 10725.9999
$ws.Range("M132").Value = -5779.16
$ws.Range("N132").Value = -15785.9999
# Row 136 (item id 43999)
$ws.Range("H136").Value = 901.4583
$ws.Range("I136").Value = 901.4583
$ws.Range("K136").Value = 2704.3749
$ws.Range("M136").Value = -154.3748999999998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94 (item id 19939)
$ws.Range("H94").Value = 11905352
$ws.Range("I94").Value = 12500519
$ws.Range("K94").Value = 12500519
$ws.Range("M94").Value = -12500068
# Row 107 (item id 27706)
$ws.Range("H107").Value = 1442.4333
$ws.Range("I107").Value = 1316.2174
$ws.Range("J107").Value = 1857.1428
$ws.Range("K107").Value = 1316.2174
$ws.Range("L107").Value = 1857.1428
$ws.Range("M107").Value = 603.7826
$ws.Range("N107").Value = -5697.1428
# Row 134 (item id 43998)
$ws.Range("H134").Value = 13929.294
$ws.Range("I134").Value = 8549.9375
$ws.Range("K134").Value = 25649.8125
$ws.Range("M134").Value = -23114.8125

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58 (item id 44021)
$ws.Range("H58").Value = 911.61536
$ws.Range("I58").Value = 911.61536
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 911.61536
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
# Row 62 (item id 12580)
$ws.Range("H62").Value = 7145796
$ws.Range("I62").Value = 3011.6667
$ws.Range("J62").Value = 50002500
$ws.Range("K62").Value = 3011.6667
$ws.Range("L62").Value = 50002500
$ws.Range("M62").Value = -2387.6667
$ws.Range("N62").Value = -50003748
# Row 65 (item id 12580)
$ws.Range("H65").Value = 7145796
$ws.Range("I65").Value = 3011.6667
$ws.Range("J65").Value = 50002500
$ws.Range("K65").Value = 15058.3335
$ws.Range("L65").Value = 250012500
$ws.Range("M65").Value = -11938.3335
$ws.Range("N65").Value = -250018740
# Row 132 (item id 44019)
$ws.Range("H132").Value = 7246.5
$ws.Range("I132").Value = 9547.357
$ws.Range("J132").Value = 3220
$ws.Range("K132").Value = 28642.071
$ws.Range("L132").Value = 9660
$ws.Range("M132").Value = -26112.071
$ws.Range("N132").Value = -14720
# Row 134 (item id 44020)
$ws.Range("H134").Value = 12346935
$ws.Range("I134").Value = 13334450
$ws.Range("K134").Value = 40003350
$ws.Range("M134").Value = -40000815
# Row 136 (item id 44021)
$ws.Range("H136").Value = 911.61536
$ws.Range("I136").Value = 911.61536
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2734.84608
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 3 (item id 44094)
$ws.Range("H3").Value = 7477.077
$ws.Range("I3").Value = 4826.9
$ws.Range("K3").Value = 14480.7
$ws.Range("M3").Value = -14368.7
# Row 68 (item id 12895)
$ws.Range("H68").Value = 1617.1395
$ws.Range("I68").Value = 655.7778
$ws.Range("J68").Value = 1871.6177
$ws.Range("K68").Value = 1967.3334
$ws.Range("L68").Value = 5614.8531
$ws.Range("M68").Value = -1156.3334
$ws.Range("N68").Value = -7236.8531
# Row 71 (item id 12895)
$ws.Range("H71").Value = 1617.1395
$ws.Range("I71").Value = 655.7778
$ws.Range("J71").Value = 1871.6177
$ws.Range("K71").Value = 5902.000199999999
$ws.Range("L71").Value = 16844.5593
$ws.Range("M71").Value = -1846.000199999999
$ws.Range("N71").Value = -24956.5593
# Row 109 (item id 27854)
$ws.Range("H109").Value = 145479.14
$ws.Range("J109").Value = 4000
$ws.Range("L109").Value = 12000
$ws.Range("N109").Value = -14080
# Row 112 (item id 27855)
$ws.Range("H112").Value = 6545.4546
$ws.Range("J112").Value = 6545.4546
$ws.Range("L112").Value = 19636.3638
$ws.Range("N112").Value = -21852.3638
# Row 131 (item id 36060)
$ws.Range("H131").Value = 21740616
$ws.Range("I131").Value = 142857820
$ws.Range("J131").Value = 1630.8718
$ws.Range("K131").Value = 428573460
$ws.Range("L131").Value = 4892.6154
$ws.Range("M131").Value = -428568420
$ws.Range("N131").Value = -14972.6154
# Row 132 (item id 43972)
$ws.Range("H132").Value = 1668
$ws.Range("I132").Value = 1004
$ws.Range("K132").Value = 9036
$ws.Range("M132").Value = -6506
# Row 140 (item id 44097)
$ws.Range("H140").Value = 30271.223
$ws.Range("I140").Value = 42321.36
$ws.Range("J140").Value = 2884.5454
$ws.Range("K140").Value = 126964.08
$ws.Range("L140").Value = 8653.636200000001
$ws.Range("M140").Value = -121784.08
$ws.Range("N140").Value = -19013.6362

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132 (item id 44008)
$ws.Range("H132").Value = 2788.476
$ws.Range("I132").Value = 2611.2307
$ws.Range("J132").Value = 3076.5
$ws.Range("K132").Value = 7833.6921
$ws.Range("L132").Value = 9229.5
$ws.Range("M132").Value = -5303.6921
$ws.Range("N132").Value = -14289.5
# Row 139 (item id 42373)
$ws.Range("H139").Value = 32000
$ws.Range("J139").Value = 32000
$ws.Range("L139").Value = 32000
$ws.Range("N139").Value = -42280

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 2 (item id 2631)
$ws.Range("H2").Value = 2647.6155
$ws.Range("J2").Value = 2536.5833
$ws.Range("L2").Value = 2536.5833
$ws.Range("N2").Value = -2760.5833
# Row 22 (item id 5277)
$ws.Range("H22").Value = 1650.125
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 1743
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1743
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -2333
# Row 27 (item id 5277)
$ws.Range("H27").Value = 1650.125
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 1743
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 1743
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -1957
# Row 132 (item id 44058)
$ws.Range("H132").Value = 54968.156
$ws.Range("I132").Value = 1949.3334
$ws.Range("J132").Value = 145857.58
$ws.Range("K132").Value = 5848.0002
$ws.Range("L132").Value = 437572.74
$ws.Range("M132").Value = -3318.0002
$ws.Range("N132").Value = -442632.74
# Row 136 (item id 44060)
$ws.Range("H136").Value = 6020.6
$ws.Range("J136").Value = 1251.875
$ws.Range("L136").Value = 3755.625
$ws.Range("N136").Value = -8855.625

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 109 (item id 27161)
$ws.Range("H109").Value = 20377
$ws.Range("J109").Value = 20377
$ws.Range("L109").Value = 20377
$ws.Range("N109").Value = -23151

